$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.708.68"
$ws.Range("E2").Value = "  +7.90%  "
Set-TextValue $ws.Range("D3") "3.534.89"
$ws.Range("E3").Value = "  +10.39%  "
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue $ws.Range("D5") "191.27"
$ws.Range("E5").Value = "  +9.84%  "
Set-TextValue $ws.Range("D6") "552.73"
$ws.Range("E6").Value = "  +7.77%  "
Set-TextValue $ws.Range("D7") "3.527.21"
$ws.Range("E7").Value = "  +10.15%  "
Set-TextValue $ws.Range("D8") "0.609"
$ws.Range("E8").Value = "  +3.61%  "
$ws.Range("E9").Value = "  -0.03%  "
Set-TextValue $ws.Range("D10") "0.640"
$ws.Range("E10").Value = "  +7.72%  "
Set-TextValue $ws.Range("D11") "56.63"
$ws.Range("E11").Value = "  +7.98%  "
$ws.Range("E12").Value = "  +17.66%  "
$ws.Range("E13").Value = "  +9.59%  "
Set-TextValue $ws.Range("D14") "9.48"
$ws.Range("E14").Value = "  +7.33%  "
Set-TextValue $ws.Range("D15") "4.092.72"
$ws.Range("E15").Value = "  +10.58%  "
Set-TextValue $ws.Range("D16") "3.533.12"
$ws.Range("E16").Value = "  +10.78%  "
Set-TextValue $ws.Range("D17") "67.740.03"
$ws.Range("E17").Value = "  +8.23%  "
Set-TextValue $ws.Range("D18") "0.122"
$ws.Range("E18").Value = "  +6.21%  "
Set-TextValue $ws.Range("D19") "18.41"
$ws.Range("E19").Value = "  +8.07%  "
Set-TextValue $ws.Range("D20") "11.91"
$ws.Range("E20").Value = "  +9.78%  "
Set-TextValue $ws.Range("D21") "0.999"
$ws.Range("E21").Value = "  +5.55%  "
Set-TextValue $ws.Range("D22") "408.54"
$ws.Range("E22").Value = "  +13.22%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D23") "4.30"
$ws.Range("E23").Value = "  +12.95%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D24") "3.95"
$ws.Range("E24").Value = "  +7.24%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "84.90"
$ws.Range("E25").Value = "  +6.72%  "
Set-TextValue $ws.Range("D26") "11.44"
$ws.Range("E26").Value = "  +3.90%  "
Set-TextValue $ws.Range("D27") "2.98"
$ws.Range("E27").Value = "  +15.40%  "
Set-TextValue $ws.Range("D28") "6.15"
$ws.Range("E28").Value = "  +2.59%  "
Set-TextValue $ws.Range("D29") "11.97"
$ws.Range("E29").Value = "  +7.48%  "
Set-TextValue $ws.Range("D30") "8.82"
$ws.Range("E30").Value = "  +9.31%  "
Set-TextValue $ws.Range("D31") "706.73"
$ws.Range("E31").Value = "  +9.33%  "
Set-TextValue $ws.Range("D32") "30.58"
$ws.Range("E32").Value = "  +9.21%  "
Set-TextValue $ws.Range("D33") "6.84"
$ws.Range("E33").Value = "  +9.97%  "
$ws.Range("E34").Value = "  +7.35%  "
$ws.Range("E35").Value = "  +8.65%  "
Set-TextValue $ws.Range("D36") "60.65"
$ws.Range("E36").Value = "  +5.96%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D37") "0.0₃0834"
$ws.Range("E37").Value = "  +21.57%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D38") "39.23"
$ws.Range("E38").Value = "  +8.40%  "
Set-TextValue $ws.Range("D39") "0.999"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  +7.58%  "
$ws.Range("E41").Value = "  +11.98%  "
Set-TextValue $ws.Range("D42") "3.37"
$ws.Range("E42").Value = "  +21.46%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D44") "2.99"
$ws.Range("E44").Value = "  +16.81%  "
Set-TextValue $ws.Range("D45") "2.70"
$ws.Range("E45").Value = "  +8.65%  "
Set-TextValue $ws.Range("D46") "3.048.57"
$ws.Range("E46").Value = "  +7.72%  "
Set-TextValue $ws.Range("D47") "3.39"
$ws.Range("E47").Value = "  +16.24%  "
Set-TextValue $ws.Range("D48") "0.0423"
$ws.Range("E48").Value = "  +10.09%  "
Set-TextValue $ws.Range("D49") "2.74"
$ws.Range("E49").Value = "  +2.65%  "
Set-TextValue $ws.Range("D50") "9.00"
$ws.Range("E50").Value = "  +20.91%  "
$ws.Range("E51").Value = "  +7.79%  "
